# Add new facilities to the "choices" sheet of client_review.xlsx
#
# The "facilities" choice list (list_name/name/label rows 17-27) is being
# expanded from 11 entries to 20 entries: several labels are simplified
# (the "Chitungwiza-"/"Marondera-"/"Chegutu-" prefixes are dropped) and nine
# brand-new facilities are appended (Goromonzi, Sanyati, Zaka, Mberengwa...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

# list_name / name / label for every row of the "facilities" choice list
$facilities = @(
    @("facilities", "seke_north",           "1. Chitungwiza-Seke North clinic"),
    @("facilities", "seke_south",           "2. Chitungwiza-Seke South clinic"),
    @("facilities", "city_med",             "3. City Med hospital"),
    @("facilities", "zengeza",              "4. Zengeza Clinic"),
    @("facilities", "chitungwiza_central",  "5. Chitungwiza Central Hospital"),
    @("facilities", "chegutu_norton",       "6. Chegutu- Norton hospital"),
    @("facilities", "chegutu_district",     "7. Chegutu District Hospital"),
    @("facilities", "monera",               "8. Monera clinic(Norton Outreach)"),
    @("facilities", "marondera",            "9. Marondera District Hospital"),
    @("facilities", "mahusekwa",            "10. Mahusekwa Hospital"),
    @("facilities", "makumbe",              "11. Goromonzi-Makumbe Mission Hospital"),
    @("facilities", "ruwa",                 "12. Goromonzi-Ruwa Rehab Hospital"),
    @("facilities", "kadoma",               "13. Sanyati-Kadoma Hospital"),
    @("facilities", "ndanga",               "14. Zaka-Ndanga District Hospital"),
    @("facilities", "musiso",               "15. Zaka-Musiso Mission Hospital"),
    @("facilities", "musiso",               "16. Mberengwa-Musiso Mission Hospital"),
    @("facilities", "musiso",               "17. Mberengwa-Mnene Mission Hospital"),
    @("facilities", "musiso",               "18. Mberengwa-Musume Mission Hospital"),
    @("facilities", "musiso",               "19. Mberengwa-Mberengwa District Hospital"),
    @("facilities", "other",                "20. Other")
)

$startRow = 17
for ($i = 0; $i -lt $facilities.Count; $i++) {
    $r = $startRow + $i
    $row = $facilities[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Re-format the rewritten rows with the plain unstyled look used elsewhere
# in the sheet (rows 17-35 lose the old orange "facilities" banner style,
# the very last row keeps a banner-like style).
$ws.Range("A17:C35").Font.Name = "Cambria"
$ws.Range("A17:C35").Font.Size = 11

# The "choices" tab becomes the active/selected sheet after this edit.
$ws.Activate()
